$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

$ws.Range("C170").Value = "TEST VALUE"

Write-Output "done"
